$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DSD")

# Row 4: REF_AREA -> GEO_PICT (rest of the row is unchanged)
$ws.Range("A4").Value = "GEO_PICT"

# Row 13 (UNIT_MEASURE): codelist renamed CL_UNIT_MEASURE -> CL_COM_UNIT_MEASURE
$ws.Range("F13").Value = "CL_COM_UNIT_MEASURE"

# Insert a new row at 14 for UNIT_MULT (pushes old rows 14-15 down to 15-16)
$ws.Rows.Item(14).Insert()
$ws.Range("A14").Value = "UNIT_MULT"
$ws.Range("B14").Value = "Unit multiplier"
$ws.Range("C14").Value = "Y"
$ws.Range("D14").Value = "Attribute"
$ws.Range("E14").Value = "Coded"
$ws.Range("F14").Value = "CL_COM_UNIT_MULT"
$ws.Range("G14").Value = "Y"

# Row 15 is now OBS_STATUS: codelist renamed CL_OBS_STATUS -> CL_COM_OBS_STATUS
$ws.Range("F15").Value = "CL_COM_OBS_STATUS"

# Insert a new row at 16 for DATA_SOURCE (pushes old COMMENT row from 15 down to 17)
$ws.Rows.Item(16).Insert()
$ws.Range("A16").Value = "DATA_SOURCE"
$ws.Range("B16").Value = "Data source"
$ws.Range("C16").Value = "N"
$ws.Range("D16").Value = "Attribute"
$ws.Range("E16").Value = "Uncoded"
$ws.Range("F16").Value = "Text"
$ws.Range("G16").Value = "N"

# Row 17 is now COMMENT: ID renamed COMMENT -> OBS_COMMENT
$ws.Range("A17").Value = "OBS_COMMENT"

# Insert a new row at 18 for CONF_STATUS
$ws.Rows.Item(18).Insert()
$ws.Range("A18").Value = "CONF_STATUS"
$ws.Range("B18").Value = "Confidentiality status"
$ws.Range("C18").Value = "Y"
$ws.Range("D18").Value = "Attribute"
$ws.Range("E18").Value = "Coded"
$ws.Range("F18").Value = "CL_COM_CONF_STATUS"
$ws.Range("G18").Value = "Y"

# Restore the selection on the DSD sheet to match the saved workbook state
$ws.Range("D23").Select()
